# Applies the row-rotation edit described in the commit diff:
# rows 179 <-> 180 swap, and rows 234/235/236/237 cyclic-rotate
# (234<-237, 235<-234, 236<-235, 237<-236). Column A (id), C (Div) and
# D (Date) are left untouched since the diff shows them unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 179
$ws.Cells.Item(179, 2).Value = 6992695
$ws.Cells.Item(179, 5).Value = "Muang Thong United"
$ws.Cells.Item(179, 6).Value = "Uthai Thani FC"
$ws.Cells.Item(179, 7).Value = 5
$ws.Cells.Item(179, 8).Value = 2
$ws.Cells.Item(179, 9).Value = 0
$ws.Cells.Item(179, 10).Value = 2
$ws.Cells.Item(179, 11).Value = "H"
$ws.Cells.Item(179, 12).Value = 2.1
$ws.Cells.Item(179, 13).Value = 3.75
$ws.Cells.Item(179, 14).Value = 2.7
$ws.Cells.Item(179, 15).Value = 1.95
$ws.Cells.Item(179, 16).Value = 3.8
$ws.Cells.Item(179, 17).Value = 2.9
$ws.Cells.Item(179, 18).Value = -0.25
$ws.Cells.Item(179, 19).Value = 1.8
$ws.Cells.Item(179, 20).Value = 2
$ws.Cells.Item(179, 21).Value = 3
$ws.Cells.Item(179, 22).Value = 1.825
$ws.Cells.Item(179, 23).Value = 1.975
$ws.Cells.Item(179, 24).Value = 0.95
$ws.Cells.Item(179, 25).Value = -1
$ws.Cells.Item(179, 26).Value = -1
$ws.Cells.Item(179, 27).Value = 0.8
$ws.Cells.Item(179, 28).Value = -1
$ws.Cells.Item(179, 29).Value = 0.825
$ws.Cells.Item(179, 30).Value = -1

# Row 180
$ws.Cells.Item(180, 2).Value = 8026714
$ws.Cells.Item(180, 5).Value = "BG Pathum United"
$ws.Cells.Item(180, 6).Value = "Buriram United"
$ws.Cells.Item(180, 7).Value = 1
$ws.Cells.Item(180, 8).Value = 1
$ws.Cells.Item(180, 9).Value = 0
$ws.Cells.Item(180, 10).Value = 1
$ws.Cells.Item(180, 11).Value = "D"
$ws.Cells.Item(180, 12).Value = 3
$ws.Cells.Item(180, 13).Value = 3.6
$ws.Cells.Item(180, 14).Value = 2
$ws.Cells.Item(180, 15).Value = 3.1
$ws.Cells.Item(180, 16).Value = 3.75
$ws.Cells.Item(180, 17).Value = 1.95
$ws.Cells.Item(180, 18).Value = 0.5
$ws.Cells.Item(180, 19).Value = 1.825
$ws.Cells.Item(180, 20).Value = 1.975
$ws.Cells.Item(180, 21).Value = 2.75
$ws.Cells.Item(180, 22).Value = 1.85
$ws.Cells.Item(180, 23).Value = 1.95
$ws.Cells.Item(180, 24).Value = -1
$ws.Cells.Item(180, 25).Value = 2.75
$ws.Cells.Item(180, 26).Value = -1
$ws.Cells.Item(180, 27).Value = 0.825
$ws.Cells.Item(180, 28).Value = -1
$ws.Cells.Item(180, 29).Value = -1
$ws.Cells.Item(180, 30).Value = 0.95

# Row 234
$ws.Cells.Item(234, 2).Value = 6992342
$ws.Cells.Item(234, 5).Value = "BG Pathum United"
$ws.Cells.Item(234, 6).Value = "Police Tero FC"
$ws.Cells.Item(234, 7).Value = 4
$ws.Cells.Item(234, 8).Value = 2
$ws.Cells.Item(234, 11).Value = "H"
$ws.Cells.Item(234, 12).Value = 1.25
$ws.Cells.Item(234, 13).Value = 6
$ws.Cells.Item(234, 14).Value = 7.5
$ws.Cells.Item(234, 15).Value = 1.25
$ws.Cells.Item(234, 16).Value = 6
$ws.Cells.Item(234, 17).Value = 7.5
$ws.Cells.Item(234, 18).Value = -2
$ws.Cells.Item(234, 19).Value = 1.95
$ws.Cells.Item(234, 20).Value = 1.85
$ws.Cells.Item(234, 21).Value = 3.75
$ws.Cells.Item(234, 22).Value = 1.85
$ws.Cells.Item(234, 23).Value = 1.95
$ws.Cells.Item(234, 24).Value = 0.25
$ws.Cells.Item(234, 25).Value = -1
$ws.Cells.Item(234, 26).Value = -1
$ws.Cells.Item(234, 27).Value = 0
$ws.Cells.Item(234, 28).Value = 0
$ws.Cells.Item(234, 29).Value = 0.8500000000000001
$ws.Cells.Item(234, 30).Value = -1

# Row 235
$ws.Cells.Item(235, 2).Value = 6992745
$ws.Cells.Item(235, 5).Value = "Sukhothai FC"
$ws.Cells.Item(235, 6).Value = "Lamphun Warrior FC"
$ws.Cells.Item(235, 7).Value = 0
$ws.Cells.Item(235, 8).Value = 3
$ws.Cells.Item(235, 11).Value = "A"
$ws.Cells.Item(235, 12).Value = 2.3
$ws.Cells.Item(235, 13).Value = 3.25
$ws.Cells.Item(235, 14).Value = 2.875
$ws.Cells.Item(235, 15).Value = 3
$ws.Cells.Item(235, 16).Value = 3.3
$ws.Cells.Item(235, 17).Value = 2.15
$ws.Cells.Item(235, 18).Value = 0.25
$ws.Cells.Item(235, 19).Value = 1.9
$ws.Cells.Item(235, 20).Value = 1.9
$ws.Cells.Item(235, 21).Value = 2.75
$ws.Cells.Item(235, 22).Value = 1.85
$ws.Cells.Item(235, 23).Value = 1.95
$ws.Cells.Item(235, 24).Value = -1
$ws.Cells.Item(235, 25).Value = -1
$ws.Cells.Item(235, 26).Value = 1.15
$ws.Cells.Item(235, 27).Value = -1
$ws.Cells.Item(235, 28).Value = 0.8999999999999999
$ws.Cells.Item(235, 29).Value = 0.425
$ws.Cells.Item(235, 30).Value = -0.5

# Row 236
$ws.Cells.Item(236, 2).Value = 6992750
$ws.Cells.Item(236, 5).Value = "Bangkok United"
$ws.Cells.Item(236, 6).Value = "Uthai Thani FC"
$ws.Cells.Item(236, 7).Value = 3
$ws.Cells.Item(236, 8).Value = 0
$ws.Cells.Item(236, 11).Value = "H"
$ws.Cells.Item(236, 12).Value = 1.3
$ws.Cells.Item(236, 13).Value = 5.5
$ws.Cells.Item(236, 14).Value = 7
$ws.Cells.Item(236, 15).Value = 1.285
$ws.Cells.Item(236, 16).Value = 6
$ws.Cells.Item(236, 17).Value = 7
$ws.Cells.Item(236, 18).Value = -1.75
$ws.Cells.Item(236, 19).Value = 1.875
$ws.Cells.Item(236, 20).Value = 1.925
$ws.Cells.Item(236, 21).Value = 3.5
$ws.Cells.Item(236, 22).Value = 1.975
$ws.Cells.Item(236, 23).Value = 1.825
$ws.Cells.Item(236, 24).Value = 0.2849999999999999
$ws.Cells.Item(236, 25).Value = -1
$ws.Cells.Item(236, 26).Value = -1
$ws.Cells.Item(236, 27).Value = 0.875
$ws.Cells.Item(236, 28).Value = -1
$ws.Cells.Item(236, 29).Value = -1
$ws.Cells.Item(236, 30).Value = 0.825

# Row 237
$ws.Cells.Item(237, 2).Value = 6992749
$ws.Cells.Item(237, 5).Value = "Buriram United"
$ws.Cells.Item(237, 6).Value = "Khonkaen United"
$ws.Cells.Item(237, 7).Value = 8
$ws.Cells.Item(237, 8).Value = 2
$ws.Cells.Item(237, 11).Value = "H"
$ws.Cells.Item(237, 12).Value = 1.166
$ws.Cells.Item(237, 13).Value = 7
$ws.Cells.Item(237, 14).Value = 11
$ws.Cells.Item(237, 15).Value = 1.166
$ws.Cells.Item(237, 16).Value = 7
$ws.Cells.Item(237, 17).Value = 10
$ws.Cells.Item(237, 18).Value = -2.25
$ws.Cells.Item(237, 19).Value = 1.9
$ws.Cells.Item(237, 20).Value = 1.9
$ws.Cells.Item(237, 21).Value = 3.75
$ws.Cells.Item(237, 22).Value = 1.95
$ws.Cells.Item(237, 23).Value = 1.85
$ws.Cells.Item(237, 24).Value = 0.1659999999999999
$ws.Cells.Item(237, 25).Value = -1
$ws.Cells.Item(237, 26).Value = -1
$ws.Cells.Item(237, 27).Value = 0.8999999999999999
$ws.Cells.Item(237, 28).Value = -1
$ws.Cells.Item(237, 29).Value = 0.95
$ws.Cells.Item(237, 30).Value = -1
